# Add data for 2022-09-15 carjacking-by-neighborhood-by-month report.
# This updates the "current month" label (through Sept 06 -> Sept 07) and
# bumps the relevant September-column monthly counts for the neighborhoods
# that had a carjacking recorded on 2022-09-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the workbook's sheet reference.
$ws.Name = "Through 2022-09-07"

# Update the "current month" column header text (column B, row 1).
$ws.Range("B1").Value = "September 2022 (through September 07)"

# Cell updates: [A1-style ref] = new value
$updates = @{
    "K3"   = 5
    "B6"   = 3
    "BM6"  = 2
    "T8"   = 1
    "AU9"  = 1
    "BD9"  = 3
    "B10"  = 1
    "BM11" = 1
    "T14"  = 2
    "AC15" = 1
    "T22"  = 1
    "AC22" = 1
    "AC23" = 1
    "K24"  = 1
    "BD36" = 1
    "AC40" = 3
    "B50"  = 2
    "AU50" = 1
    "B55"  = 1
    "AU98" = 2
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
